# NON-COVID Player Stats / K / K_aggregate.xlsx
# Rebuild the kicker FG% table: add a "Season Group" column and expand
# each player into Group1 / Group2 / Difference rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Season Group" column between Player and FG% ---
$ws.Columns("B:B").Insert()

# --- 2. Expand each existing single player row into three rows
#        (Group1 / Group2 / Difference), preserving the row's banding style ---
# Existing rows (after the column insert) are at 2,3,4,5,6 (one per player).
# Insert two new rows directly below each player row, starting with the
# bottom-most player row first so earlier row numbers stay valid, and each
# new row inherits the style of the row immediately above it.
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").Insert()
$ws.Rows("6:6").Insert()
$ws.Rows("6:6").Insert()
$ws.Rows("5:5").Insert()
$ws.Rows("5:5").Insert()
$ws.Rows("4:4").Insert()
$ws.Rows("4:4").Insert()
$ws.Rows("3:3").Insert()
$ws.Rows("3:3").Insert()

# --- 3. Header row ---
$ws.Cells.Item(1,1).Value = "Player"
$ws.Cells.Item(1,2).Value = "Season Group"
$ws.Cells.Item(1,3).Value = "FG%"
$ws.Cells.Item(1,4).Value = "Lng"

# --- 4. Data rows ---
$ws.Cells.Item(2,1).Value = "Jason Myers"
$ws.Cells.Item(2,2).Value = "Group1"
$ws.Cells.Item(2,3).Value = 85.33333333333333
$ws.Cells.Item(2,4).Value = 56

$ws.Cells.Item(3,1).Value = "Jason Myers"
$ws.Cells.Item(3,2).Value = "Group2"
$ws.Cells.Item(3,3).Value = 87.3
$ws.Cells.Item(3,4).Value = 56.66666666666666

$ws.Cells.Item(4,1).Value = "Jason Myers"
$ws.Cells.Item(4,2).Value = "Difference"
$ws.Cells.Item(4,3).Value = 1.966666666666669
$ws.Cells.Item(4,4).Value = 0.6666666666666643

$ws.Cells.Item(5,1).Value = "Justin Tucker"
$ws.Cells.Item(5,2).Value = "Group1"
$ws.Cells.Item(5,3).Value = 93.63333333333333
$ws.Cells.Item(5,4).Value = 57.33333333333334

$ws.Cells.Item(6,1).Value = "Justin Tucker"
$ws.Cells.Item(6,2).Value = "Group2"
$ws.Cells.Item(6,3).Value = 81.93333333333334
$ws.Cells.Item(6,4).Value = 54.66666666666666

$ws.Cells.Item(7,1).Value = "Justin Tucker"
$ws.Cells.Item(7,2).Value = "Difference"
$ws.Cells.Item(7,3).Value = -11.69999999999999
$ws.Cells.Item(7,4).Value = -2.666666666666671

$ws.Cells.Item(8,1).Value = "Matt Gay"
$ws.Cells.Item(8,2).Value = "Group1"
$ws.Cells.Item(8,3).Value = 86.23333333333333
$ws.Cells.Item(8,4).Value = 54.66666666666666

$ws.Cells.Item(9,1).Value = "Matt Gay"
$ws.Cells.Item(9,2).Value = "Group2"
$ws.Cells.Item(9,3).Value = 85.86666666666667
$ws.Cells.Item(9,4).Value = 57

$ws.Cells.Item(10,1).Value = "Matt Gay"
$ws.Cells.Item(10,2).Value = "Difference"
$ws.Cells.Item(10,3).Value = -0.36666666666666
$ws.Cells.Item(10,4).Value = 2.333333333333336

$ws.Cells.Item(11,1).Value = "sportsref download (73)"
$ws.Cells.Item(11,2).Value = "Group1"
$ws.Cells.Item(11,3).Value = 85.33333333333333
$ws.Cells.Item(11,4).Value = 56

$ws.Cells.Item(12,1).Value = "sportsref download (73)"
$ws.Cells.Item(12,2).Value = "Group2"
$ws.Cells.Item(12,3).Value = 87.3
$ws.Cells.Item(12,4).Value = 56.66666666666666

$ws.Cells.Item(13,1).Value = "sportsref download (73)"
$ws.Cells.Item(13,2).Value = "Difference"
$ws.Cells.Item(13,3).Value = 1.966666666666669
$ws.Cells.Item(13,4).Value = 0.6666666666666643

$ws.Cells.Item(14,1).Value = "Younghoe Koo"
$ws.Cells.Item(14,2).Value = "Group1"
$ws.Cells.Item(14,3).Value = 92.16666666666667
$ws.Cells.Item(14,4).Value = 52

$ws.Cells.Item(15,1).Value = "Younghoe Koo"
$ws.Cells.Item(15,2).Value = "Group2"
$ws.Cells.Item(15,3).Value = 82.16666666666667
$ws.Cells.Item(15,4).Value = 55.33333333333334

$ws.Cells.Item(16,1).Value = "Younghoe Koo"
$ws.Cells.Item(16,2).Value = "Difference"
$ws.Cells.Item(16,3).Value = -10
$ws.Cells.Item(16,4).Value = 3.333333333333336

# --- 5. Reset page margins to Excel defaults (as in the saved file) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
